$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 261, pushing the existing rows 261-275 down to 262-276.
$ws.Rows.Item(261).Insert()

# Populate the newly inserted row 261 with the new record's data.
$ws.Range("A261").Value = 7
$ws.Range("B261").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C261").Value = "Ñuble"
$ws.Range("D261").Value = 44585
$ws.Range("E261").Value = 16
$ws.Range("F261").Value = 100114001
$ws.Range("G261").Value = "Papa"
$ws.Range("H261").Value = "Patagonia"
$ws.Range("I261").Value = "1a nueva(o)"
$ws.Range("J261").Value = 200
$ws.Range("K261").Value = 7000
$ws.Range("L261").Value = 7500
$ws.Range("M261").Value = 7250
$ws.Range("N261").Value = "`$/saco 25 kilos"
$ws.Range("O261").Value = "Provincia de Diguillín"
$ws.Range("P261").Value = 290
$ws.Range("Q261").Value = 25
$ws.Range("R261").Value = "Hortaliza"
